# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Rebuild DeliveryPlan / VehicleLog / TruckUsageLog sheets with corrected data.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) DeliveryPlan sheet: 16 data rows -> 14 data rows (17 total -> 15 total)
# -----------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("DeliveryPlan")

# Remove the two trailing rows that no longer exist (old rows 16 & 17)
$wsPlan.Rows.Item(17).Delete() | Out-Null
$wsPlan.Rows.Item(16).Delete() | Out-Null

$planData = @(
    @("20240105-PLANT_001-DC_001-LARGE-#1","MAT_A|PLANT_001|DC_001|2024-01-06|net demand for forecast|000041","MAT_A","PLANT_001","DC_001",45297,45296,45298,12,"LARGE",0.996,0.996,0.9825),
    @("20240105-PLANT_001-DC_001-LARGE-#1","MAT_A|PLANT_001|DC_001|2024-01-06|net demand for forecast|000089","MAT_A","PLANT_001","DC_001",45297,45296,45298,33,"LARGE",0.996,0.996,0.9825),
    @("20240105-PLANT_001-DC_001-LARGE-#1","MAT_B|PLANT_001|DC_001|2024-01-06|net demand for forecast|000094","MAT_B","PLANT_001","DC_001",45297,45296,45298,16,"LARGE",0.996,0.996,0.9825),
    @("20240105-PLANT_001-DC_001-LARGE-#1","MAT_A|PLANT_001|DC_001|2024-01-07|net demand for forecast|000090","MAT_A","PLANT_001","DC_001",45298,45296,45298,33,"LARGE",0.996,0.996,0.9825),
    @("20240105-PLANT_001-DC_001-LARGE-#1","MAT_B|PLANT_001|DC_001|2024-01-07|net demand for forecast|000095","MAT_B","PLANT_001","DC_001",45298,45296,45298,11,"LARGE",0.996,0.996,0.9825),
    @("20240105-PLANT_001-DC_001-LARGE-#2","MAT_B|PLANT_001|DC_001|2024-01-07|net demand for forecast|000095","MAT_B","PLANT_001","DC_001",45298,45296,45298,5,"LARGE",0.956,0.956,0.9375),
    @("20240105-PLANT_001-DC_001-LARGE-#2","MAT_A|PLANT_001|DC_001|2024-01-08|net demand for forecast|000091","MAT_A","PLANT_001","DC_001",45299,45296,45298,33,"LARGE",0.956,0.956,0.9375),
    @("20240105-PLANT_001-DC_001-LARGE-#2","MAT_B|PLANT_001|DC_001|2024-01-08|net demand for forecast|000096","MAT_B","PLANT_001","DC_001",45299,45296,45298,16,"LARGE",0.956,0.956,0.9375),
    @("20240105-PLANT_001-DC_001-LARGE-#2","MAT_A|PLANT_001|DC_001|2024-01-09|net demand for forecast|000092","MAT_A","PLANT_001","DC_001",45300,45296,45298,33,"LARGE",0.956,0.956,0.9375),
    @("20240105-PLANT_001-DC_001-LARGE-#2","MAT_B|PLANT_001|DC_001|2024-01-09|net demand for forecast|000097","MAT_B","PLANT_001","DC_001",45300,45296,45299,16,"LARGE",0.956,0.956,0.9375),
    @("20240105-PLANT_001-DC_002-MEDIUM-#1","MAT_A|PLANT_001|DC_002|2024-01-06|net demand for forecast|000037","MAT_A","PLANT_001","DC_002",45297,45296,45298,44,"MEDIUM",1,1,1),
    @("20240105-PLANT_001-DC_002-MEDIUM-#1","MAT_A|PLANT_001|DC_002|2024-01-06|net demand for forecast|000085","MAT_A","PLANT_001","DC_002",45297,45296,45298,16,"MEDIUM",1,1,1),
    @("20240105-PLANT_001-DC_002-MEDIUM-#2","MAT_A|PLANT_001|DC_002|2024-01-06|net demand for forecast|000085","MAT_A","PLANT_001","DC_002",45297,45296,45299,51,"MEDIUM",1,1,1),
    @("20240105-PLANT_001-DC_002-MEDIUM-#2","MAT_A|PLANT_001|DC_002|2024-01-07|net demand for forecast|000086","MAT_A","PLANT_001","DC_002",45298,45296,45298,9,"MEDIUM",1,1,1)
)

$rowCount = $planData.Count
$colCount = 13
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r,$c] = $planData[$r][$c]
    }
}
$wsPlan.Range("A2:M15").Value = $arr

# -----------------------------------------------------------------
# 2) VehicleLog sheet: 3 data rows -> 4 data rows
# -----------------------------------------------------------------
$wsVeh = $wb.Worksheets.Item("VehicleLog")

# insert a new row so there is room for the extra LARGE-#2 vehicle
$wsVeh.Rows.Item(3).Insert() | Out-Null

$vehData = @(
    @(45296,"PLANT_001","DC_001","LARGE",1,"20240105-PLANT_001-DC_001-LARGE-#1",105,99.59999999999999,196.5,0.996,0.9825,"threshold"),
    @(45296,"PLANT_001","DC_001","LARGE",2,"20240105-PLANT_001-DC_001-LARGE-#2",103,95.59999999999999,187.5,0.956,0.9375,"threshold"),
    @(45296,"PLANT_001","DC_002","MEDIUM",1,"20240105-PLANT_001-DC_002-MEDIUM-#1",60,60,120,1,1,"threshold"),
    @(45296,"PLANT_001","DC_002","MEDIUM",2,"20240105-PLANT_001-DC_002-MEDIUM-#2",60,60,120,1,1,"threshold")
)

$rowCount = $vehData.Count
$colCount = 12
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r,$c] = $vehData[$r][$c]
    }
}
$wsVeh.Range("A2:L5").Value = $arr
$wsVeh.Range("A2:A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# -----------------------------------------------------------------
# 3) TruckUsageLog sheet: truck_used for PLANT_001/DC_001/LARGE goes 1 -> 2
# -----------------------------------------------------------------
$wsTruck = $wb.Worksheets.Item("TruckUsageLog")
$wsTruck.Cells.Item(2,5).Value = 2

# -----------------------------------------------------------------
# 4) ValidationLog sheet: content unchanged (only internal string table
#    ordering differed upstream) - nothing to do.
# -----------------------------------------------------------------
